# "Generate Report for Archive"
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet
#    that reports it (Overview!E2:F2, zh-cn!C2, de-de!C2 - all share the
#    same "Status" string).
# 2) The "Status" column got narrower (report column resized) on all
#    three sheets: Overview columns E & F, and column C on the zh-cn /
#    de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the Status values ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Resize the Status columns ---
# Target stored column width is 13.4101845877511 characters; Excel's
# ColumnWidth property only accepts/reports width in whole-pixel
# increments, so feed it the character width (12.5) that rounds to the
# closest representable stored width.
$newColumnWidth = 12.5

$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
